$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 11-13 take on the "colored" style (same as row 2 / row 8 originally: s=2,3,4)
$ws.Range("A2:C2").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$ws.Range("A2:C2").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Range("A2:C2").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)

# Row 14 takes on the "continuation" style (same as row 3: s=5,6,7)
$ws.Range("A3:C3").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)

# Row 8 changes from "colored" style to "continuation" style
$ws.Range("A3:C3").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)

# Now fill in the values
$ws.Range("A11").Value = 41888
$ws.Range("B11").Value = "13.10"
$ws.Range("C11").Value = "More backbacks dropped in houses or guaranteed backpack in ""Wildness"" / ""Houses"" type as loot"

$ws.Range("A12").Value = 41888
$ws.Range("B12").Value = "13.10"
$ws.Range("C12").Value = "Loot Vehicles should have empty cargo"

$ws.Range("A13").Value = 41888
$ws.Range("B13").Value = "13.10"
$ws.Range("C13").Value = "Disable red markers in occupied houses"

$ws.Range("A14").Value = 41888
$ws.Range("B14").Value = "13.10"
$ws.Range("C14").Value = "More scopes dropped"

# Clear the selection marker (set current selection to A1)
$ws.Range("A1").Select()
